# Update InsideBet Data: Automatizado
# Applies updated stat values to the RESUMEN_STATS_Bundesliga worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - Hamburger SV
$ws.Range("C8").Value = 26.6
$ws.Range("D8").Value = 47.5
$ws.Range("E8").Value = 22
$ws.Range("F8").Value = 242
$ws.Range("G8").Value = 1980
$ws.Range("H8").Value = 22
$ws.Range("I8").Value = 25
$ws.Range("K8").Value = 44
$ws.Range("L8").Value = 24
$ws.Range("O8").Value = 54
$ws.Range("R8").Value = 0.86
$ws.Range("S8").Value = 2
$ws.Range("T8").Value = 1.09
$ws.Range("U8").Value = 1.95

# Row 13 - Mainz 05
$ws.Range("E13").Value = 23
$ws.Range("F13").Value = 253
$ws.Range("G13").Value = 2070
$ws.Range("H13").Value = 23
$ws.Range("I13").Value = 26
$ws.Range("J13").Value = 15
$ws.Range("K13").Value = 41
$ws.Range("L13").Value = 19
$ws.Range("O13").Value = 56
$ws.Range("Q13").Value = 1.13
$ws.Range("R13").Value = 0.65
$ws.Range("S13").Value = 1.78
$ws.Range("T13").Value = 0.83
$ws.Range("U13").Value = 1.48

# Row 14 - RB Leipzig
$ws.Range("C14").Value = 26

# Row 16 - Stuttgart
$ws.Range("C16").Value = 26
